# "complete the ood_ness test" -- fill in the final simulated mean/std
# results (rows 2-31, i.e. trial index 0-29) on the sole worksheet. Rows
# 2-21 already existed with placeholder numbers that get overwritten with
# final values; rows 22-31 are brand new trials (index 20-29) appended
# below the existing data.
$data = @(
    @(2, 0, 0.5, 0),
    @(3, 1, 0.712500023841858, 0.08164963862597714),
    @(4, 2, 0.712500023841858, 0.05803495229609949),
    @(5, 3, 0.8833333492279053, 0.08186201268405847),
    @(6, 4, 0.5541666865348815, 0.04289846085055259),
    @(7, 5, 0.6708333492279053, 0.0533593689555652),
    @(8, 6, 0.8916666746139527, 0.1869900315084763),
    @(9, 7, 0.5541666746139526, 0.08079467945472155),
    @(10, 8, 0.9458333611488342, 0.02825970591949954),
    @(11, 9, 0.5166666984558106, 0.07500000728501242),
    @(12, 10, 0.4568070933915337, 0.04687115721599745),
    @(13, 11, 0.702345260690466, 0.08932050717143956),
    @(14, 12, 0.7038370508766112, 0.06233519582119591),
    @(15, 13, 0.8830210268383762, 0.08195720948512833),
    @(16, 14, 0.5133808155440034, 0.06630544977308069),
    @(17, 15, 0.6587419753435316, 0.060854347344914),
    @(18, 16, 0.8689760537784854, 0.2320513957851537),
    @(19, 17, 0.5075267970623552, 0.1069716999002428),
    @(20, 18, 0.9457443970496, 0.02840381171191628),
    @(21, 19, 0.512892120145961, 0.07680164309621712),
    @(22, 20, 0.4961805555555556, 0.006892164319888624),
    @(23, 21, 0.784375, 0.1101871303905484),
    @(24, 22, 0.7911458333333332, 0.05794451269038275),
    @(25, 23, 0.9635416666666666, 0.03377171226607729),
    @(26, 24, 0.6088541666666667, 0.0702934216229143),
    @(27, 25, 0.7392361111111112, 0.07728639419972352),
    @(28, 26, 0.923611111111111, 0.1476245812028222),
    @(29, 27, 0.6371527777777779, 0.02131955754584547),
    @(30, 28, 0.9927083333333334, 0.005851492898039144),
    @(31, 29, 0.5368055555555556, 0.07994711552560352)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
}

# New rows 22-31 need column A formatted like the existing data rows (bold,
# centered, bordered "s=1" style) -- copy the format from row 21's A cell.
$ws.Range("A21").Copy()
$ws.Range("A22:A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
